# Update the LINE_TRIALS_URL sheet: fill in row 2 of the trial-tracking
# table with the first trial entry, then move the selection to I3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LINE_TRIALS_URL")
$ws.Activate()

$ws.Range("A2").Value = "Vendor A"
$ws.Range("B2").Value = "SEALANT (POTTING)"
$ws.Range("C2").Value = "Completed"

# E2 needs the same date number format as D2 already carries; copy that
# formatting over before writing the value.
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("D2").Value = 36892
$ws.Range("E2").Value = 36923
$ws.Range("F2").Value = "ok"
$ws.Range("G2").Value = "ok"

$ws.Range("I3").Select()
